# Add Betting Markets Analytics (Handicap, O/U, BTTS)
#
# The underlying source data for the "Olympique Lyonnais" stats sheet got
# re-ordered: the player that used to occupy row 11 (Rachid Ghezzal) and the
# player that used to occupy row 12 (Pavel Šulc) swap places. Columns A
# (League) and B (Team) are identical for both rows and stay put; every
# other column from C (Player_Name) through DK (type) needs to be exchanged
# between the two rows (column DL, goalsPrevented, is blank in both rows and
# is left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowA = 11
$rowB = 12
$firstCol = "C"
$lastCol  = "DK"
# (DL, goalsPrevented, is blank on both rows both before and after, so it is
# intentionally left out of the swap to avoid disturbing that empty cell.)

$rangeA = $ws.Range("$firstCol${rowA}:$lastCol${rowA}")
$rangeB = $ws.Range("$firstCol${rowB}:$lastCol${rowB}")

# Grab both rows worth of data (as 2D arrays), then write each back into the
# other row so the two players' full stat lines trade places.
$valuesA = $rangeA.Value2
$valuesB = $rangeB.Value2

$rangeA.Value2 = $valuesB
$rangeB.Value2 = $valuesA
